$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3891, 3891, 4036, 4036, 4055, 4145, 4402, 4402, 4461, 4480, 4480)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
